$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Regenerated "K" column (column G) values - save_data was regenerated to
# compute K directly instead of deriving it from Strike#, so the stored
# strike-count values change for every row (std/mean and s_vals were
# recalculated upstream from this new K column).
$kValues = @(
    0,  # G2
    5,  # G3
    2,  # G4
    2,  # G5
    3,  # G6
    5,  # G7
    2,  # G8
    4,  # G9
    3,  # G10
    5,  # G11
    4,  # G12
    2,  # G13
    8,  # G14
    5,  # G15
    3,  # G16
    5,  # G17
    1,  # G18
    1,  # G19
    5,  # G20
    8,  # G21
    5,  # G22
    5,  # G23
    9,  # G24
    2,  # G25
    4,  # G26
    1,  # G27
    6,  # G28
    3,  # G29
    1,  # G30
    2,  # G31
    3,  # G32
    0,  # G33
    0   # G34
)

$startRow = 2
for ($i = 0; $i -lt $kValues.Count; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 7).Value = $kValues[$i]
}
